$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-25 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2026-02-26 Thursday", 2)

# Update the multiplication table cells by explicit (row, column) coordinates
# to avoid any ambiguity from duplicate/overlapping values.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "70×87=6090"
$t.Cell(1,2).Range.Text  = "23×59=1357"
$t.Cell(1,3).Range.Text  = "92×62=5704"
$t.Cell(1,4).Range.Text  = "90×16=1440"
$t.Cell(1,5).Range.Text  = "49×62=3038"

$t.Cell(5,1).Range.Text  = "93×73=6789"
$t.Cell(5,2).Range.Text  = "87×86=7482"
$t.Cell(5,3).Range.Text  = "22×82=1804"
$t.Cell(5,4).Range.Text  = "74×90=6660"
$t.Cell(5,5).Range.Text  = "86×59=5074"

$t.Cell(10,1).Range.Text = "17×69=1173"
$t.Cell(10,2).Range.Text = "12×50=600"
$t.Cell(10,3).Range.Text = "54×73=3942"
$t.Cell(10,4).Range.Text = "23×81=1863"
$t.Cell(10,5).Range.Text = "74×74=5476"

$t.Cell(15,1).Range.Text = "50×24=1200"
$t.Cell(15,2).Range.Text = "77×16=1232"
$t.Cell(15,3).Range.Text = "74×74=5476"
$t.Cell(15,4).Range.Text = "39×42=1638"
$t.Cell(15,5).Range.Text = "19×54=1026"

$t.Cell(20,1).Range.Text = "57×54=3078"
$t.Cell(20,2).Range.Text = "70×46=3220"
$t.Cell(20,3).Range.Text = "41×22=902"
$t.Cell(20,4).Range.Text = "84×75=6300"
$t.Cell(20,5).Range.Text = "47×74=3478"
